$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1: update "Last status check on" timestamp text
$ws.Range("F1").Value = "Last status check on: 15.02.2022 14:30"

# Row 7 updates
$ws.Range("B7").Value = 37.9
$ws.Range("C7").Value = 37.5

# D7: must become a literal text string "+0.4" (not auto-converted to a number),
# and end up with the default (no explicit) cell style.
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "+0.4"
$ws.Range("D7").ClearFormats()

# E7: must become a literal text string (not parsed into a date/number),
# and lose its previous date-style formatting.
$ws.Range("E7").ClearFormats()
$ws.Range("E7").Value = "2022-02-15 14:34:58"
